$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name and card number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long card-number-like string; pre-format as Text so Excel
# doesn't coerce it into a double (which would lose the literal digits /
# render in scientific notation).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 03.12.2024"

# Row 6 (existing transaction, update dates/desc/amount)
$ws.Range("B6").Value = "04.12."
$ws.Range("C6").Value = "05.12."
$ws.Range("D6").Value = "KARTENZ./04.12 LIDL RO"
$ws.Range("E6").Value = "58,42-"

# Row 7
$ws.Range("B7").Value = "08.12."
$ws.Range("C7").Value = "09.12."
$ws.Range("D7").Value = "ZALANDO MKTPLC EU DHYANA"
$ws.Range("E7").Value = "165,90-"

# Row 8
$ws.Range("B8").Value = "12.12."
$ws.Range("C8").Value = "13.12."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 26013271"
$ws.Range("E8").Value = "84,68-"

# Row 9
$ws.Range("B9").Value = "14.12."
$ws.Range("C9").Value = "15.12."
$ws.Range("D9").Value = "PAYPAL TBFUNN"
$ws.Range("E9").Value = "68,15-"

# Row 10 (was a blank filler row, now gets a real transaction)
$ws.Range("B10").Value = "16.12."
$ws.Range("C10").Value = "17.12."
$ws.Range("D10").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E10").Value = "25,03-"
# Amount column switches from the "blank filler" style (wrap + vcenter) to
# the normal right-aligned amount style used by the other rows.
$ws.Range("E10").WrapText = $false
$ws.Range("E10").VerticalAlignment = -4107

# Row 11 (was a blank filler row, now gets a real transaction)
$ws.Range("B11").Value = "20.12."
$ws.Range("C11").Value = "21.12."
$ws.Range("D11").Value = "EBAY MKTPLC EU UDNLJK"
$ws.Range("E11").Value = "101,29-"
$ws.Range("E11").WrapText = $false
$ws.Range("E11").VerticalAlignment = -4107

# Row 12: closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 22.12.2024"
$ws.Range("E12").Value = "503,47-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 01.01.2025"
